$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Formula = "=""9.25***"""
$ws.Range("C4").Formula = "=""181.9***"""
$ws.Range("D4").Formula = "=""7.68***"""
$ws.Range("B5").Formula = "=""(2.04)"""
$ws.Range("C5").Formula = "=""(50.8)"""
$ws.Range("B6").Formula = "=""0.060"""
$ws.Range("C6").Formula = "=""120.7"""
$ws.Range("D6").Formula = "=""38.0*"""
$ws.Range("B7").Formula = "=""(21.1)"""
$ws.Range("C7").Formula = "=""(528.6)"""
$ws.Range("B8").Formula = "=""10.3***"""
$ws.Range("C8").Formula = "=""189.1***"""
$ws.Range("B9").Formula = "=""(2.46)"""
$ws.Range("C9").Formula = "=""(50.9)"""
$ws.Range("D9").Formula = "=""(2.42)"""
$ws.Range("B10").Formula = "=""-56.6***"""
$ws.Range("C10").Formula = "=""-941.3***"""
$ws.Range("D10").Formula = "=""56.5***"""
$ws.Range("B11").Formula = "=""(1.42)"""
$ws.Range("C11").Formula = "=""(27.3)"""
$ws.Range("B12").Formula = "=""-47.4***"""
$ws.Range("C12").Formula = "=""-759.4***"""
$ws.Range("B13").Formula = "=""(1.46)"""
$ws.Range("C13").Formula = "=""(42.9)"""
$ws.Range("D13").Formula = "=""(1.83)"""
$ws.Range("B14").Formula = "=""-10.3"""
$ws.Range("C14").Formula = "=""-68.5"""
$ws.Range("D14").Formula = "=""33.9"""
$ws.Range("B15").Formula = "=""(22.7)"""
$ws.Range("C15").Formula = "=""(556.6)"""
$ws.Range("B16").Formula = "=""14.9"""
$ws.Range("C16").Formula = "=""280.3"""
$ws.Range("D16").Formula = "=""-39.8*"""
$ws.Range("B17").Formula = "=""(22.1)"""
$ws.Range("C17").Formula = "=""(551.6)"""
$ws.Range("B18").Formula = "=""4.58"""
$ws.Range("C18").Formula = "=""211.9***"""
$ws.Range("B19").Formula = "=""(3.55)"""
$ws.Range("C19").Formula = "=""(59.5)"""
$ws.Range("B22").Formula = "=""0.65"""
$ws.Range("C22").Formula = "=""0.90"""
$ws.Range("D22").Formula = "=""0.13"""
$ws.Range("B23").Formula = "=""0.65"""
$ws.Range("C23").Formula = "=""0.90"""
$ws.Range("B24").Formula = "=""0.65"""
$ws.Range("C24").Formula = "=""0.90"""
$ws.Range("B25").Formula = "=""0.67"""
$ws.Range("C25").Formula = "=""0.55"""
$ws.Range("D25").Formula = "=""0.068"""
